# Architecture slide (slide 7): fix the "Node/Debian Dig Microservice" info
# box - nudge its vertical position back up (auto-reflow side effect of the
# text change below) and add the missing word "HTTP" to the end of the
# "Service for performing dig requests using" line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Locate the textbox by its stable shape Id (27 in the OOXML / p:cNvPr id)
# rather than a hard-coded collection index, in case shape ordering ever
# shifts.
$targetId = 27
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq $targetId) {
        $shape = $candidate
    }
}

# --- 1) Reposition the shape -------------------------------------------
# Target OOXML offset is y="2346376" EMU (x is untouched). Shape.Top is in
# points (1 pt = 12700 EMU) and is rounded to a 32-bit float internally by
# the host, so feed it the closest value that round-trips to exactly the
# right EMU count instead of the naive 2346376/12700 division.
$shape.Top = 184.7540283203125

# --- 2) Add the missing word "HTTP" -------------------------------------
# InsertAfter() on the shape's full TextRange appends at the very end of
# the text (end of the last paragraph) and inherits the run formatting
# already in effect there (sz=1600, accent1), which is exactly what the
# new run needs, so no further formatting calls are required.
$tr = $shape.TextFrame.TextRange
$newRun = $tr.InsertAfter(" HTTP")
